{"js": "// Change \"Version 2.\" -> \"Version 1.\" (wireframes version bump revert).\n// Use search() so the edit merges the split runs the same way Word would\n// when you retype across an existing run boundary.\n\n// \"Versi\" + \"on\" are two separate runs in the source; search finds the\n// logical text spanning them and replacing it collapses them into one run.\nconst versionResults = context.document.body.search(\"Version\", { matchCase: true });\nversionResults.load(\"text\");\nawait context.sync();\n\nif (versionResults.items.length > 0) {\n  versionResults.items[0].insertText(\"Version\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// \" 2\" and the trailing \".\" are two separate runs around a bookmark;\n// replacing the whole \"2.\" span merges them into a single \" 1.\" run and\n// drops the now-empty trailing run, moving the bookmark to the end.\nconst numberResults = context.document.body.search(\"2.\", { matchCase: true });\nnumberResults.load(\"text\");\nawait context.sync();\n\nif (numberResults.items.length > 0) {\n  numberResults.items[0].insertText(\"1.\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Change \"Version 2.\" -> \"Version 1.\" (wireframes version bump revert).\n$d = $word.ActiveDocument\n\n# \"Versi\" + \"on\" are two separate runs in the source; a Find/Replace over\n# the full word collapses them into a single \"Version\" run, same as\n# retyping across the run boundary in Word.\n$d.Content.Find.Execute(\"Version\", $false, $false, $false, $false, $false, $true, 1, $false, \"Version\", 2)\n\n# Flip the digit only (leave the trailing \".\" run and the \"_GoBack\"\n# bookmark between them untouched for now) so the bookmark position is\n# preserved instead of being swallowed by a wider Find/Replace.\n$d.Content.Find.Execute(\"2\", $false, $false, $false, $false, $false, $true, 1, $false, \"1\", 2)\n\n# Now fold the lone \".\" run into the \" 1\" run, the way Word does when the\n# trailing run collapses: delete the \".\" character, then re-insert it\n# immediately before the (now-trailing) \"_GoBack\" bookmark so it merges\n# into the preceding run and the bookmark ends up after it.\n$full = $d.Content\n$dotRange = $d.Range($full.End - 2, $full.End - 1)\n$dotRange.Delete()\n\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$insertPoint = $d.Range($bm.Start, $bm.Start)\n$insertPoint.InsertBefore(\".\")\n"}
